# Adds rows 14-33 to the General Log sheet, logging subsequent compound
# upload attempts (duplicates + successful submissions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Cells.Item(14, 1).Value = '2024-09-25 11:12:50'
$ws.Cells.Item(14, 2).Value = 'ORM-0515839'
$ws.Cells.Item(14, 3).Value = 'Z195631098'
$ws.Cells.Item(14, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(14, 5).Value = 'Duplicate'
$ws.Cells.Item(14, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'

# Row 15
$ws.Cells.Item(15, 1).Value = '2024-09-25 11:12:51'
$ws.Cells.Item(15, 2).Value = 'ORM-0515840'
$ws.Cells.Item(15, 3).Value = 'Z2754556176'
$ws.Cells.Item(15, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(15, 5).Value = 'Duplicate'
$ws.Cells.Item(15, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'

# Row 16
$ws.Cells.Item(16, 1).Value = '2024-09-25 11:13:29'
$ws.Cells.Item(16, 2).Value = 'ORM-0515839'
$ws.Cells.Item(16, 3).Value = 'Z195631098'
$ws.Cells.Item(16, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(16, 5).Value = 'Duplicate'
$ws.Cells.Item(16, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'

# Row 17
$ws.Cells.Item(17, 1).Value = '2024-09-25 11:13:31'
$ws.Cells.Item(17, 2).Value = 'ORM-0515840'
$ws.Cells.Item(17, 3).Value = 'Z2754556176'
$ws.Cells.Item(17, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(17, 5).Value = 'Duplicate'
$ws.Cells.Item(17, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'

# Row 18
$ws.Cells.Item(18, 1).Value = '2024-09-25 11:15:05'
$ws.Cells.Item(18, 2).Value = 'ORM-0515839'
$ws.Cells.Item(18, 3).Value = 'Z195631098'
$ws.Cells.Item(18, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(18, 5).Value = 'Duplicate'
$ws.Cells.Item(18, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'

# Row 19
$ws.Cells.Item(19, 1).Value = '2024-09-25 11:15:07'
$ws.Cells.Item(19, 2).Value = 'ORM-0515840'
$ws.Cells.Item(19, 3).Value = 'Z2754556176'
$ws.Cells.Item(19, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(19, 5).Value = 'Duplicate'
$ws.Cells.Item(19, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'

# Row 20
$ws.Cells.Item(20, 1).Value = '2024-09-25 11:16:21'
$ws.Cells.Item(20, 2).Value = 'ORM-0515848'
$ws.Cells.Item(20, 3).Value = 'Z195631098'
$ws.Cells.Item(20, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(20, 5).Value = 'Success'
$ws.Cells.Item(20, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(20, 7).Value = 201
$ws.Cells.Item(20, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c6d5b10e4653161b72f1"},"data":{"type":"material","id":"asset:66f3c6d5b10e4653161b72f1","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c6d5b10e4653161b72f1"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3c6d5b10e4653161b72f1","id":"asset:66f3c6d5b10e4653161b72f1","eid":"asset:66f3c6d5b10e4653161b72f1","name":"ORM-0515848","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-09-25T08:16:21.164Z","editedAt":"2024-09-25T08:16:21.164Z","type":"asset","digest":"38834240","fields":{"Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515848"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3c6d5b10e4653161b72f2","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c6d5b10e4653161b72f2"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c6d5b10e4653161b72f1/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3c6d5b10e4653161b72f1"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3c6d5b10e4653161b72f2","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c6d5b10e4653161b72f2"},"attributes":{"type":"batch","eid":"batch:66f3c6d5b10e4653161b72f2","name":"ORM-0515848-001","digest":"48077144","fields":{"Batch Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515848-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:16:21.480Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"78459056","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}},{"type":"materialDrawing","id":"asset:66f3c6d5b10e4653161b72f1","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c6d5b10e4653161b72f1/drawing?format=cdxml"},"attributes":{"id":"asset:66f3c6d5b10e4653161b72f1","type":"CHEMICAL_DRAWING"}}]}'

# Row 21
$ws.Cells.Item(21, 1).Value = '2024-09-25 11:16:24'
$ws.Cells.Item(21, 2).Value = 'ORM-0515849'
$ws.Cells.Item(21, 3).Value = 'Z2754556176'
$ws.Cells.Item(21, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(21, 5).Value = 'Success'
$ws.Cells.Item(21, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(21, 7).Value = 201
$ws.Cells.Item(21, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c6d7a01482132101b337"},"data":{"type":"material","id":"asset:66f3c6d7a01482132101b337","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c6d7a01482132101b337"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3c6d7a01482132101b337","id":"asset:66f3c6d7a01482132101b337","eid":"asset:66f3c6d7a01482132101b337","name":"ORM-0515849","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-09-25T08:16:23.917Z","editedAt":"2024-09-25T08:16:23.917Z","type":"asset","digest":"39131108","fields":{"Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0515849"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3c6d8a01482132101b338","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c6d8a01482132101b338"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c6d7a01482132101b337/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3c6d7a01482132101b337"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3c6d8a01482132101b338","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c6d8a01482132101b338"},"attributes":{"type":"batch","eid":"batch:66f3c6d8a01482132101b338","name":"ORM-0515849-001","digest":"86043289","fields":{"Batch Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide, hydrogen bromide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;&amp;middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Name":{"value":"ORM-0515849-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:16:24.220Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"78459056","fields":{}}},{"type":"materialDrawing","id":"asset:66f3c6d7a01482132101b337","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c6d7a01482132101b337/drawing?format=cdxml"},"attributes":{"id":"asset:66f3c6d7a01482132101b337","type":"CHEMICAL_DRAWING"}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 22
$ws.Cells.Item(22, 1).Value = '2024-09-25 11:22:12'
$ws.Cells.Item(22, 2).Value = 'ORM-0515850'
$ws.Cells.Item(22, 3).Value = 'Z195631098'
$ws.Cells.Item(22, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(22, 5).Value = 'Success'
$ws.Cells.Item(22, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(22, 7).Value = 201
$ws.Cells.Item(22, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c833a01482132101b347"},"data":{"type":"material","id":"asset:66f3c833a01482132101b347","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c833a01482132101b347"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3c833a01482132101b347","id":"asset:66f3c833a01482132101b347","eid":"asset:66f3c833a01482132101b347","name":"ORM-0515850","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-09-25T08:22:11.916Z","editedAt":"2024-09-25T08:22:11.916Z","type":"asset","digest":"51238311","fields":{"Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515850"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3c834a01482132101b348","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c834a01482132101b348"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c833a01482132101b347/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3c833a01482132101b347"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"materialDrawing","id":"asset:66f3c833a01482132101b347","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c833a01482132101b347/drawing?format=cdxml"},"attributes":{"id":"asset:66f3c833a01482132101b347","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"batch:66f3c834a01482132101b348","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c834a01482132101b348"},"attributes":{"type":"batch","eid":"batch:66f3c834a01482132101b348","name":"ORM-0515850-001","digest":"16964069","fields":{"Batch Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515850-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:22:12.229Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"61526003","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 23
$ws.Cells.Item(23, 1).Value = '2024-09-25 11:22:15'
$ws.Cells.Item(23, 2).Value = 'ORM-0515851'
$ws.Cells.Item(23, 3).Value = 'Z2754556176'
$ws.Cells.Item(23, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(23, 5).Value = 'Success'
$ws.Cells.Item(23, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(23, 7).Value = 201
$ws.Cells.Item(23, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c837399ff04b78d3b950"},"data":{"type":"material","id":"asset:66f3c837399ff04b78d3b950","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c837399ff04b78d3b950"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3c837399ff04b78d3b950","id":"asset:66f3c837399ff04b78d3b950","eid":"asset:66f3c837399ff04b78d3b950","name":"ORM-0515851","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-09-25T08:22:15.127Z","editedAt":"2024-09-25T08:22:15.127Z","type":"asset","digest":"76371280","fields":{"Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0515851"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3c837399ff04b78d3b951","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c837399ff04b78d3b951"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c837399ff04b78d3b950/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3c837399ff04b78d3b950"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3c837399ff04b78d3b951","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c837399ff04b78d3b951"},"attributes":{"type":"batch","eid":"batch:66f3c837399ff04b78d3b951","name":"ORM-0515851-001","digest":"65911264","fields":{"Batch Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide, hydrogen bromide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;&amp;middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Name":{"value":"ORM-0515851-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:22:15.310Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"materialDrawing","id":"asset:66f3c837399ff04b78d3b950","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c837399ff04b78d3b950/drawing?format=cdxml"},"attributes":{"id":"asset:66f3c837399ff04b78d3b950","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"61526003","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 24
$ws.Cells.Item(24, 1).Value = '2024-09-25 11:25:15'
$ws.Cells.Item(24, 2).Value = 'ORM-0515852'
$ws.Cells.Item(24, 3).Value = 'Z195631098'
$ws.Cells.Item(24, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(24, 5).Value = 'Success'
$ws.Cells.Item(24, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(24, 7).Value = 201
$ws.Cells.Item(24, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c8eaccf303112a9a6319"},"data":{"type":"material","id":"asset:66f3c8eaccf303112a9a6319","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c8eaccf303112a9a6319"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3c8eaccf303112a9a6319","id":"asset:66f3c8eaccf303112a9a6319","eid":"asset:66f3c8eaccf303112a9a6319","name":"ORM-0515852","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-09-25T08:25:14.449Z","editedAt":"2024-09-25T08:25:14.449Z","type":"asset","digest":"50162357","fields":{"Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515852"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3c8eaccf303112a9a631a","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c8eaccf303112a9a631a"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c8eaccf303112a9a6319/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3c8eaccf303112a9a6319"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3c8eaccf303112a9a631a","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c8eaccf303112a9a631a"},"attributes":{"type":"batch","eid":"batch:66f3c8eaccf303112a9a631a","name":"ORM-0515852-001","digest":"94979605","fields":{"Batch Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515852-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:25:14.924Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"48247138","fields":{}}},{"type":"materialDrawing","id":"asset:66f3c8eaccf303112a9a6319","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c8eaccf303112a9a6319/drawing?format=cdxml"},"attributes":{"id":"asset:66f3c8eaccf303112a9a6319","type":"CHEMICAL_DRAWING"}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 25
$ws.Cells.Item(25, 1).Value = '2024-09-25 11:25:18'
$ws.Cells.Item(25, 2).Value = 'ORM-0515853'
$ws.Cells.Item(25, 3).Value = 'Z2754556176'
$ws.Cells.Item(25, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(25, 5).Value = 'Success'
$ws.Cells.Item(25, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(25, 7).Value = 201
$ws.Cells.Item(25, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c8ed8ac1f915a4ada772"},"data":{"type":"material","id":"asset:66f3c8ed8ac1f915a4ada772","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c8ed8ac1f915a4ada772"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3c8ed8ac1f915a4ada772","id":"asset:66f3c8ed8ac1f915a4ada772","eid":"asset:66f3c8ed8ac1f915a4ada772","name":"ORM-0515853","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-09-25T08:25:17.401Z","editedAt":"2024-09-25T08:25:17.401Z","type":"asset","digest":"27686644","fields":{"Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0515853"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3c8ed8ac1f915a4ada773","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c8ed8ac1f915a4ada773"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c8ed8ac1f915a4ada772/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3c8ed8ac1f915a4ada772"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3c8ed8ac1f915a4ada773","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c8ed8ac1f915a4ada773"},"attributes":{"type":"batch","eid":"batch:66f3c8ed8ac1f915a4ada773","name":"ORM-0515853-001","digest":"90612978","fields":{"Batch Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide, hydrogen bromide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;&amp;middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Name":{"value":"ORM-0515853-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:25:17.893Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"materialDrawing","id":"asset:66f3c8ed8ac1f915a4ada772","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c8ed8ac1f915a4ada772/drawing?format=cdxml"},"attributes":{"id":"asset:66f3c8ed8ac1f915a4ada772","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"48247138","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 26
$ws.Cells.Item(26, 1).Value = '2024-09-25 11:29:52'
$ws.Cells.Item(26, 2).Value = 'ORM-0515854'
$ws.Cells.Item(26, 3).Value = 'Z195631098'
$ws.Cells.Item(26, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(26, 5).Value = 'Success'
$ws.Cells.Item(26, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(26, 7).Value = 201
$ws.Cells.Item(26, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c9ff399ff04b78d3b952"},"data":{"type":"material","id":"asset:66f3c9ff399ff04b78d3b952","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c9ff399ff04b78d3b952"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3c9ff399ff04b78d3b952","id":"asset:66f3c9ff399ff04b78d3b952","eid":"asset:66f3c9ff399ff04b78d3b952","name":"ORM-0515854","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-09-25T08:29:51.259Z","editedAt":"2024-09-25T08:29:51.259Z","type":"asset","digest":"61891900","fields":{"Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515854"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3c9ff399ff04b78d3b953","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c9ff399ff04b78d3b953"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c9ff399ff04b78d3b952/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3c9ff399ff04b78d3b952"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3c9ff399ff04b78d3b953","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3c9ff399ff04b78d3b953"},"attributes":{"type":"batch","eid":"batch:66f3c9ff399ff04b78d3b953","name":"ORM-0515854-001","digest":"74433293","fields":{"Batch Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515854-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:29:51.916Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"26314572","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}},{"type":"materialDrawing","id":"asset:66f3c9ff399ff04b78d3b952","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3c9ff399ff04b78d3b952/drawing?format=cdxml"},"attributes":{"id":"asset:66f3c9ff399ff04b78d3b952","type":"CHEMICAL_DRAWING"}}]}'

# Row 27
$ws.Cells.Item(27, 1).Value = '2024-09-25 11:29:55'
$ws.Cells.Item(27, 2).Value = 'ORM-0515855'
$ws.Cells.Item(27, 3).Value = 'Z2754556176'
$ws.Cells.Item(27, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(27, 5).Value = 'Success'
$ws.Cells.Item(27, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(27, 7).Value = 201
$ws.Cells.Item(27, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca02399ff04b78d3b954"},"data":{"type":"material","id":"asset:66f3ca02399ff04b78d3b954","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca02399ff04b78d3b954"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3ca02399ff04b78d3b954","id":"asset:66f3ca02399ff04b78d3b954","eid":"asset:66f3ca02399ff04b78d3b954","name":"ORM-0515855","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-09-25T08:29:54.458Z","editedAt":"2024-09-25T08:29:54.458Z","type":"asset","digest":"83080449","fields":{"Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0515855"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3ca02399ff04b78d3b955","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3ca02399ff04b78d3b955"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca02399ff04b78d3b954/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3ca02399ff04b78d3b954"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"materialDrawing","id":"asset:66f3ca02399ff04b78d3b954","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca02399ff04b78d3b954/drawing?format=cdxml"},"attributes":{"id":"asset:66f3ca02399ff04b78d3b954","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"batch:66f3ca02399ff04b78d3b955","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3ca02399ff04b78d3b955"},"attributes":{"type":"batch","eid":"batch:66f3ca02399ff04b78d3b955","name":"ORM-0515855-001","digest":"99076106","fields":{"Batch Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide, hydrogen bromide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;&amp;middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Name":{"value":"ORM-0515855-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:29:54.923Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"26314572","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 28
$ws.Cells.Item(28, 1).Value = '2024-09-25 11:31:32'
$ws.Cells.Item(28, 2).Value = 'ORM-0515856'
$ws.Cells.Item(28, 3).Value = 'Z195631098'
$ws.Cells.Item(28, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(28, 5).Value = 'Success'
$ws.Cells.Item(28, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(28, 7).Value = 201
$ws.Cells.Item(28, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca63399ff04b78d3b956"},"data":{"type":"material","id":"asset:66f3ca63399ff04b78d3b956","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca63399ff04b78d3b956"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3ca63399ff04b78d3b956","id":"asset:66f3ca63399ff04b78d3b956","eid":"asset:66f3ca63399ff04b78d3b956","name":"ORM-0515856","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-09-25T08:31:31.876Z","editedAt":"2024-09-25T08:31:31.876Z","type":"asset","digest":"31175883","fields":{"Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515856"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3ca64399ff04b78d3b957","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3ca64399ff04b78d3b957"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca63399ff04b78d3b956/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3ca63399ff04b78d3b956"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3ca64399ff04b78d3b957","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3ca64399ff04b78d3b957"},"attributes":{"type":"batch","eid":"batch:66f3ca64399ff04b78d3b957","name":"ORM-0515856-001","digest":"53712521","fields":{"Batch Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515856-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:31:32.093Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"materialDrawing","id":"asset:66f3ca63399ff04b78d3b956","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca63399ff04b78d3b956/drawing?format=cdxml"},"attributes":{"id":"asset:66f3ca63399ff04b78d3b956","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"27915298","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 29
$ws.Cells.Item(29, 1).Value = '2024-09-25 11:31:35'
$ws.Cells.Item(29, 2).Value = 'ORM-0515857'
$ws.Cells.Item(29, 3).Value = 'Z2754556176'
$ws.Cells.Item(29, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(29, 5).Value = 'Success'
$ws.Cells.Item(29, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(29, 7).Value = 201
$ws.Cells.Item(29, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca66b10e4653161b7309"},"data":{"type":"material","id":"asset:66f3ca66b10e4653161b7309","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca66b10e4653161b7309"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3ca66b10e4653161b7309","id":"asset:66f3ca66b10e4653161b7309","eid":"asset:66f3ca66b10e4653161b7309","name":"ORM-0515857","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-09-25T08:31:34.474Z","editedAt":"2024-09-25T08:31:34.474Z","type":"asset","digest":"40461342","fields":{"Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0515857"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3ca66b10e4653161b730a","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3ca66b10e4653161b730a"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca66b10e4653161b7309/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3ca66b10e4653161b7309"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3ca66b10e4653161b730a","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3ca66b10e4653161b730a"},"attributes":{"type":"batch","eid":"batch:66f3ca66b10e4653161b730a","name":"ORM-0515857-001","digest":"70452462","fields":{"Batch Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide, hydrogen bromide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;&amp;middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Name":{"value":"ORM-0515857-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:31:34.908Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"materialDrawing","id":"asset:66f3ca66b10e4653161b7309","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3ca66b10e4653161b7309/drawing?format=cdxml"},"attributes":{"id":"asset:66f3ca66b10e4653161b7309","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"27915298","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 30
$ws.Cells.Item(30, 1).Value = '2024-09-25 11:35:05'
$ws.Cells.Item(30, 2).Value = 'ORM-0515858'
$ws.Cells.Item(30, 3).Value = 'Z195631098'
$ws.Cells.Item(30, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(30, 5).Value = 'Success'
$ws.Cells.Item(30, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(30, 7).Value = 201
$ws.Cells.Item(30, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cb38399ff04b78d3b95a"},"data":{"type":"material","id":"asset:66f3cb38399ff04b78d3b95a","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cb38399ff04b78d3b95a"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3cb38399ff04b78d3b95a","id":"asset:66f3cb38399ff04b78d3b95a","eid":"asset:66f3cb38399ff04b78d3b95a","name":"ORM-0515858","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-09-25T08:35:04.906Z","editedAt":"2024-09-25T08:35:04.906Z","type":"asset","digest":"76671961","fields":{"Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515858"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3cb39399ff04b78d3b95b","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3cb39399ff04b78d3b95b"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cb38399ff04b78d3b95a/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3cb38399ff04b78d3b95a"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3cb39399ff04b78d3b95b","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3cb39399ff04b78d3b95b"},"attributes":{"type":"batch","eid":"batch:66f3cb39399ff04b78d3b95b","name":"ORM-0515858-001","digest":"25909751","fields":{"Batch Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515858-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:35:05.166Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"25407789","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}},{"type":"materialDrawing","id":"asset:66f3cb38399ff04b78d3b95a","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cb38399ff04b78d3b95a/drawing?format=cdxml"},"attributes":{"id":"asset:66f3cb38399ff04b78d3b95a","type":"CHEMICAL_DRAWING"}}]}'

# Row 31
$ws.Cells.Item(31, 1).Value = '2024-09-25 11:35:08'
$ws.Cells.Item(31, 2).Value = 'ORM-0515859'
$ws.Cells.Item(31, 3).Value = 'Z2754556176'
$ws.Cells.Item(31, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(31, 5).Value = 'Success'
$ws.Cells.Item(31, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(31, 7).Value = 201
$ws.Cells.Item(31, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cb3b399ff04b78d3b95c"},"data":{"type":"material","id":"asset:66f3cb3b399ff04b78d3b95c","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cb3b399ff04b78d3b95c"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3cb3b399ff04b78d3b95c","id":"asset:66f3cb3b399ff04b78d3b95c","eid":"asset:66f3cb3b399ff04b78d3b95c","name":"ORM-0515859","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-09-25T08:35:07.478Z","editedAt":"2024-09-25T08:35:07.478Z","type":"asset","digest":"12743630","fields":{"Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0515859"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3cb3b399ff04b78d3b95d","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3cb3b399ff04b78d3b95d"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cb3b399ff04b78d3b95c/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3cb3b399ff04b78d3b95c"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"materialDrawing","id":"asset:66f3cb3b399ff04b78d3b95c","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cb3b399ff04b78d3b95c/drawing?format=cdxml"},"attributes":{"id":"asset:66f3cb3b399ff04b78d3b95c","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"batch:66f3cb3b399ff04b78d3b95d","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3cb3b399ff04b78d3b95d"},"attributes":{"type":"batch","eid":"batch:66f3cb3b399ff04b78d3b95d","name":"ORM-0515859-001","digest":"43866544","fields":{"Batch Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide, hydrogen bromide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;&amp;middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Name":{"value":"ORM-0515859-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:35:07.911Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"25407789","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 32
$ws.Cells.Item(32, 1).Value = '2024-09-25 11:37:29'
$ws.Cells.Item(32, 2).Value = 'ORM-0515860'
$ws.Cells.Item(32, 3).Value = 'Z195631098'
$ws.Cells.Item(32, 4).Value = 'C18H19CL2N5OS'
$ws.Cells.Item(32, 5).Value = 'Success'
$ws.Cells.Item(32, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(32, 7).Value = 201
$ws.Cells.Item(32, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cbc8b10e4653161b730b"},"data":{"type":"material","id":"asset:66f3cbc8b10e4653161b730b","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cbc8b10e4653161b730b"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3cbc8b10e4653161b730b","id":"asset:66f3cbc8b10e4653161b730b","eid":"asset:66f3cbc8b10e4653161b730b","name":"ORM-0515860","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-09-25T08:37:28.263Z","editedAt":"2024-09-25T08:37:28.263Z","type":"asset","digest":"74629343","fields":{"Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515860"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3cbc8b10e4653161b730c","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3cbc8b10e4653161b730c"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cbc8b10e4653161b730b/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3cbc8b10e4653161b730b"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"materialDrawing","id":"asset:66f3cbc8b10e4653161b730b","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cbc8b10e4653161b730b/drawing?format=cdxml"},"attributes":{"id":"asset:66f3cbc8b10e4653161b730b","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"batch:66f3cbc8b10e4653161b730c","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3cbc8b10e4653161b730c"},"attributes":{"type":"batch","eid":"batch:66f3cbc8b10e4653161b730c","name":"ORM-0515860-001","digest":"24931436","fields":{"Batch Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;18&lt;/sub&gt;H&lt;sub&gt;19&lt;/sub&gt;Cl&lt;sub&gt;2&lt;/sub&gt;N&lt;sub&gt;5&lt;/sub&gt;OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0515860-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:37:28.909Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"97254194","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 33
$ws.Cells.Item(33, 1).Value = '2024-09-25 11:37:32'
$ws.Cells.Item(33, 2).Value = 'ORM-0515861'
$ws.Cells.Item(33, 3).Value = 'Z2754556176'
$ws.Cells.Item(33, 4).Value = 'C17H28N4O2'
$ws.Cells.Item(33, 5).Value = 'Success'
$ws.Cells.Item(33, 6).Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Cells.Item(33, 7).Value = 201
$ws.Cells.Item(33, 8).Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cbcbccf303112a9a6325"},"data":{"type":"material","id":"asset:66f3cbcbccf303112a9a6325","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cbcbccf303112a9a6325"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"66f3cbcbccf303112a9a6325","id":"asset:66f3cbcbccf303112a9a6325","eid":"asset:66f3cbcbccf303112a9a6325","name":"ORM-0515861","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-09-25T08:37:31.457Z","editedAt":"2024-09-25T08:37:31.457Z","type":"asset","digest":"51579799","fields":{"Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0515861"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:66f3cbcbccf303112a9a6326","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3cbcbccf303112a9a6326"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cbcbccf303112a9a6325/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:66f3cbcbccf303112a9a6325"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:66f3cbcbccf303112a9a6326","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:66f3cbcbccf303112a9a6326"},"attributes":{"type":"batch","eid":"batch:66f3cbcbccf303112a9a6326","name":"ORM-0515861-001","digest":"54944541","fields":{"Batch Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide, hydrogen bromide"},"Batch Molecular Formula":{"value":"C&lt;sub&gt;17&lt;/sub&gt;H&lt;sub&gt;28&lt;/sub&gt;N&lt;sub&gt;4&lt;/sub&gt;O&lt;sub&gt;2&lt;/sub&gt;&amp;middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Name":{"value":"ORM-0515861-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-09-25T08:37:31.903Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"materialDrawing","id":"asset:66f3cbcbccf303112a9a6325","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:66f3cbcbccf303112a9a6325/drawing?format=cdxml"},"attributes":{"id":"asset:66f3cbcbccf303112a9a6325","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"97254194","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

